# repull data, push all data, mean calculation
# Update column F ("dSF") values to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -1
$ws.Range("F8").Value = -7
$ws.Range("F9").Value = -8
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = 5
$ws.Range("F18").Value = 4
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 3
$ws.Range("F24").Value = -6
$ws.Range("F25").Value = 2
$ws.Range("F26").Value = 4
